$wb = $excel.ActiveWorkbook

# "展览" and "全部类型" sheets both list the same exhibition events.
# Bump the "想去人数" (want-to-go count) for two events:
#   row 2 (LZ 栋子动漫游戏嘉年华): 9 -> 10
#   row 4 (龙泉ACG动漫游戏博览会): 1455 -> 1456
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 10
    $ws.Range("F4").Value = 1456
}
